$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "63.368.34"
$ws.Range("E2").Value = "  +2.26%  "
$ws.Range("D3").Value = "3.475.85"
$ws.Range("E3").Value = "  +1.54%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "581.76"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").Value = "147.21"
$ws.Range("E6").Value = "  +1.57%  "
$ws.Range("D7").Value = "3.475.45"
$ws.Range("E7").Value = "  +1.53%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "0.477"
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("E10").Value = "  +1.19%  "
$ws.Range("E11").Value = "  +1.15%  "
$ws.Range("D12").Value = "0.404"
$ws.Range("E12").Value = "  +4.87%  "
$ws.Range("D13").Value = "4.071.58"
$ws.Range("E13").Value = "  +1.60%  "
$ws.Range("D14").Value = "29.64"
$ws.Range("E14").Value = "  +4.71%  "
$ws.Range("E15").Value = "  +2.30%  "
$ws.Range("D16").Value = "3.478.81"
$ws.Range("E16").Value = "  +1.61%  "
$ws.Range("E17").Value = "  +1.14%  "
$ws.Range("D18").Value = "63.412.84"
$ws.Range("E18").Value = "  +2.27%  "
$ws.Range("D19").Value = "6.37"
$ws.Range("E19").Value = "  +3.17%  "
$ws.Range("D20").Value = "14.46"
$ws.Range("E20").Value = "  +3.33%  "
$ws.Range("D21").Value = "9.34"
$ws.Range("D22").Value = "389.77"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").Value = "0.565"
$ws.Range("E23").Value = "  +2.27%  "
$ws.Range("D24").Value = "74.88"
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "3.622.10"
$ws.Range("E26").Value = "  +1.73%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").Value = "  -4.46%  "
$ws.Range("D29").Value = "7.63"
$ws.Range("E29").Value = "  +2.17%  "
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("D31").Value = "8.25"
$ws.Range("E31").Value = "  +2.48%  "
$ws.Range("D34").Value = "1.37"
$ws.Range("E34").Value = "  -3.56%  "
$ws.Range("D35").Value = "23.54"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("E36").Value = "  +0.98%  "
$ws.Range("E37").Value = "  +2.81%  "
$ws.Range("E38").Value = "  +9.34%  "
$ws.Range("D39").Value = "31.81"
$ws.Range("E39").Value = "  +12.26%  "
$ws.Range("D40").Value = "169.30"
$ws.Range("E40").Value = "  +0.79%  "
$ws.Range("D41").Value = "3.515.13"
$ws.Range("E41").Value = "  +1.70%  "
$ws.Range("D42").Value = "0.0764"
$ws.Range("E42").Value = "  +1.24%  "
$ws.Range("D43").Value = "0.800"
$ws.Range("E43").Value = "  +1.66%  "
$ws.Range("E44").Value = "  +3.64%  "
$ws.Range("D45").Value = "42.39"
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("E46").Value = "  +2.81%  "
$ws.Range("D47").Value = "4.41"
$ws.Range("E47").Value = "  -0.81%  "
$ws.Range("D48").Value = "2.603.99"
$ws.Range("E48").Value = "  +3.02%  "
$ws.Range("E49").Value = "  +9.83%  "
$ws.Range("D50").Value = "23.13"
$ws.Range("E50").Value = "  +0.89%  "
$ws.Range("E51").Value = "  +2.85%  "
